# Convolutional network size calculator.xlsx - update the MNIST sheet's
# first conv-layer input size (A33 / A35) from 512 to 128. Every other
# numeric change in the target diff is a cached formula result that
# recalculates automatically off of these two input cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MNIST")
[void]$ws.Activate()

$ws.Range("A33").Value = 128
$ws.Range("A35").Value = 128

# Reflect the author's final cursor position/selection on the MNIST sheet.
[void]$ws.Range("Z49").Select()

Write-Output "done"
